$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression
$ws.Range("B2").Value = 3128315078485088
$ws.Range("C2").Value = 3128315078485089
$ws.Range("D2").Value = 3128315078485089

# Row 3 - RandomForestRegressor
$ws.Range("B3").Value = 9305033831737.779
$ws.Range("C3").Value = 5788009779664.292
$ws.Range("D3").Value = 4380131092985.438

# Row 4 - GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 2505915285617.425
$ws.Range("C4").Value = 2505915285617.425
$ws.Range("D4").Value = 2427605432941.882

# Row 5 - AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 96834530933064.5
$ws.Range("C5").Value = 176904018845061.4
$ws.Range("D5").Value = 185001067515417.7
